# Update RF (column I) values for rows 16 through 62 from 11.25911111111111
# to 15.15071428571428, per "Update of 2025 data and RF changes".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 15.15071428571428

for ($row = 16; $row -le 62; $row++) {
    $ws.Cells.Item($row, 9).Value = $newValue  # Column I = 9
}
